$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2: year headers, extend from Q2 (2030) to AK2 (2050) with years 2031..2050
$startCol = 18  # R
$startYear = 2031
for ($i = 0; $i -lt 20; $i++) {
    $col = $startCol + $i
    $year = $startYear + $i
    $ws.Cells.Item(2, $col).Value = $year
}

# Row 3: demand values for R3:AK3
$row3vals = @(
    251444.83518433457,
    258783.94924618702,
    265263.30896215944,
    273117.41240514588,
    276039.79001293099,
    283755.32674828876,
    292067.76979976828,
    300823.30140487646,
    309563.90928671957,
    318261.39967951499,
    327702.95572908991,
    337102.51639970695,
    346536.78194036102,
    354315.41854397382,
    361902.60212297749,
    371535.66604170779,
    380982.5927209584,
    390256.83753794763,
    398846.22348578909,
    407185.65063201322
)
for ($i = 0; $i -lt $row3vals.Length; $i++) {
    $col = $startCol + $i
    $ws.Cells.Item(3, $col).Value = $row3vals[$i]
}

# Row 7: extend shared formula H7:Q7 ("=X2") to H7:AK7
for ($col = 18; $col -le 37; $col++) {
    $colLetter = $ws.Cells.Item(2, $col).Address($false, $false) -replace '\d', ''
    $ws.Cells.Item(7, $col).Formula = "=" + $colLetter + "2"
}

# Row 8: extend shared formula H8:Q8 ("=X3*3.6/1000") to H8:AK8
for ($col = 18; $col -le 37; $col++) {
    $colLetter = $ws.Cells.Item(3, $col).Address($false, $false) -replace '\d', ''
    $ws.Cells.Item(8, $col).Formula = "=" + $colLetter + "3*3.6/1000"
    $ws.Cells.Item(8, $col).NumberFormat = "0.00"
}

# Update dimension / selection via UsedRange recalculation will happen automatically.
$ws.Range("Q8:AK8").Select()

$wb.Save()
